# Auto-generated cell value updates for Spriggan_Profits workbook
# (scheduled runner refresh of market-board derived leve profit figures)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1185.7142
$ws.Range("I11").Value = 1185.7142
$ws.Range("K11").Value = 1185.7142
$ws.Range("M11").Value = -1045.7142
$ws.Range("H12").Value = 6799.6665
$ws.Range("I12").Value = 6799.6665
$ws.Range("K12").Value = 6799.6665
$ws.Range("M12").Value = -6629.6665
$ws.Range("H33").Value = 842.1667
$ws.Range("I33").Value = 625.5625
$ws.Range("K33").Value = 625.5625
$ws.Range("M33").Value = -396.5625
$ws.Range("H86").Value = 3499.5
$ws.Range("I86").Value = 3499.5
$ws.Range("K86").Value = 3499.5
$ws.Range("M86").Value = -2376.5
$ws.Range("H89").Value = 3499.5
$ws.Range("I89").Value = 3499.5
$ws.Range("K89").Value = 17497.5
$ws.Range("M89").Value = -11881.5
$ws.Range("H92").Value = 628.26666
$ws.Range("I92").Value = 833.55554
$ws.Range("K92").Value = 833.55554
$ws.Range("M92").Value = 414.44446
$ws.Range("H112").Value = 34217.53
$ws.Range("I112").Value = 2492.6667
$ws.Range("J112").Value = 62210.06
$ws.Range("K112").Value = 7478.000100000001
$ws.Range("L112").Value = 186630.18
$ws.Range("M112").Value = -6370.000100000001
$ws.Range("N112").Value = -188846.18
$ws.Range("H131").Value = 1421.6666
$ws.Range("I131").Value = 1421.6666
$ws.Range("K131").Value = 4264.9998
$ws.Range("M131").Value = 775.0002000000004

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4592.5
$ws.Range("J63").Value = 4592.5
$ws.Range("L63").Value = 4592.5
$ws.Range("N63").Value = -5964.5
$ws.Range("H66").Value = 4592.5
$ws.Range("J66").Value = 4592.5
$ws.Range("L66").Value = 22962.5
$ws.Range("N66").Value = -29826.5
$ws.Range("H74").Value = 34487520
$ws.Range("I74").Value = 47624980
$ws.Range("J74").Value = 1687.5
$ws.Range("K74").Value = 47624980
$ws.Range("L74").Value = 1687.5
$ws.Range("M74").Value = -47624106
$ws.Range("N74").Value = -3435.5
$ws.Range("H77").Value = 34487520
$ws.Range("I77").Value = 47624980
$ws.Range("J77").Value = 1687.5
$ws.Range("K77").Value = 238124900
$ws.Range("L77").Value = 8437.5
$ws.Range("M77").Value = -238120532
$ws.Range("N77").Value = -17173.5
$ws.Range("H97").Value = 403
$ws.Range("I97").Value = 403
$ws.Range("K97").Value = 403
$ws.Range("M97").Value = 93
$ws.Range("H102").Value = 8929507
$ws.Range("I102").Value = 8929507
$ws.Range("K102").Value = 8929507
$ws.Range("M102").Value = -8927885
$ws.Range("H132").Value = 3033824
$ws.Range("I132").Value = 3033824
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9101472
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9098942
$ws.Range("N132").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 985.125
$ws.Range("I94").Value = 1014
$ws.Range("J94").Value = 898.5
$ws.Range("K94").Value = 1014
$ws.Range("L94").Value = 898.5
$ws.Range("M94").Value = -563
$ws.Range("N94").Value = -1800.5
$ws.Range("H99").Value = 1824.3793
$ws.Range("I99").Value = 913.9091
$ws.Range("J99").Value = 2380.7778
$ws.Range("K99").Value = 913.9091
$ws.Range("L99").Value = 2380.7778
$ws.Range("M99").Value = 584.0909
$ws.Range("N99").Value = -5376.7778
$ws.Range("H107").Value = 49451.19
$ws.Range("I107").Value = 1343.3846
$ws.Range("K107").Value = 1343.3846
$ws.Range("M107").Value = 576.6153999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10730.842
$ws.Range("J31").Value = 13094.333
$ws.Range("L31").Value = 13094.333
$ws.Range("N31").Value = -13684.333
$ws.Range("H34").Value = 10730.842
$ws.Range("J34").Value = 13094.333
$ws.Range("L34").Value = 13094.333
$ws.Range("N34").Value = -13498.333
$ws.Range("H99").Value = 1674.125
$ws.Range("J99").Value = 1771.1666
$ws.Range("L99").Value = 1771.1666
$ws.Range("N99").Value = -4767.1666
$ws.Range("H126").Value = 1674.125
$ws.Range("J126").Value = 1771.1666
$ws.Range("L126").Value = 5313.4998
$ws.Range("N126").Value = -10253.4998
$ws.Range("H127").Value = 110000
$ws.Range("J127").Value = 110000
$ws.Range("L127").Value = 110000
$ws.Range("N127").Value = -119920

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 278.05884
$ws.Range("J2").Value = 302.1
$ws.Range("L2").Value = 1812.6
$ws.Range("N2").Value = -2038.6
$ws.Range("H62").Value = 2450
$ws.Range("J62").Value = 2400
$ws.Range("L62").Value = 7200
$ws.Range("N62").Value = -8572
$ws.Range("H65").Value = 2450
$ws.Range("J65").Value = 2400
$ws.Range("L65").Value = 21600
$ws.Range("N65").Value = -28464
$ws.Range("H98").Value = 1290
$ws.Range("I98").Value = 1290
$ws.Range("K98").Value = 3870
$ws.Range("M98").Value = -2372
$ws.Range("H107").Value = 1316.8462
$ws.Range("J107").Value = 1813.4706
$ws.Range("L107").Value = 5440.4118
$ws.Range("N107").Value = -9280.4118
$ws.Range("H117").Value = 1792.4
$ws.Range("J117").Value = 1976.6923
$ws.Range("L117").Value = 5930.0769
$ws.Range("N117").Value = -12814.0769

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 29990
$ws.Range("J47").Value = 29990
$ws.Range("L47").Value = 29990
$ws.Range("N47").Value = -31126

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2730.5186
$ws.Range("J16").Value = 4255.6665
$ws.Range("L16").Value = 4255.6665
$ws.Range("N16").Value = -4595.6665
$ws.Range("H46").Value = 695.8182
$ws.Range("I46").Value = 679.3333
$ws.Range("J46").Value = 770
$ws.Range("K46").Value = 679.3333
$ws.Range("L46").Value = 770
$ws.Range("M46").Value = -491.3333
$ws.Range("N46").Value = -1146
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2983.5676
$ws.Range("I96").Value = 1644
$ws.Range("J96").Value = 3296.1333
$ws.Range("K96").Value = 1644
$ws.Range("L96").Value = 3296.1333
$ws.Range("M96").Value = -271
$ws.Range("N96").Value = -6042.1333
$ws.Range("H135").Value = 74362.44500000001
$ws.Range("J135").Value = 74362.44500000001
$ws.Range("L135").Value = 74362.44500000001
$ws.Range("N135").Value = -84502.44500000001
